$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (Row 1) ---
$ws.Range("A1").Value = "TÍTULO"
$ws.Range("B1").Value = "DESCRIPCIÓN"
$ws.Range("C1").Value = "REFERENCIA"
$ws.Range("D1").Value = "PRECIO"
$ws.Range("E1").Value = "ISBN"
$ws.Range("F1").Value = "SECCIÓN"
$ws.Range("G1").Value = "ESTADO"
$ws.Range("H1").Value = "DESCRIPCIÓN DEL ESTADO"
$ws.Range("I1").Value = "OPERACIÓN"
$ws.Range("J1").Value = "STOCK"
$ws.Range("K1").Value = "FECHA DE PUBLICACIÓN"
$ws.Range("L1").Value = "FORMA DE ENVÍO"
$ws.Range("M1").Value = "GASTOS FIJOS"

# Copy header style from an existing styled header cell (A1) to the new header cells
$ws.Range("A1").Copy()
$ws.Range("J1:M1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

$shipNote = "Envíos muy rápidos con tarifa plana, combine discos y pague solo por el primer lote."

# Columns that must stay as TEXT even though they look numeric
$ws.Range("C2:C4").NumberFormat = "@"
$ws.Range("F2:G4").NumberFormat = "@"
$ws.Range("M2:M4").NumberFormat = "@"
# D2, E2, D4, E4 need to remain present as empty (but existing) text cells
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("E4").NumberFormat = "@"

# --- Row 2: Carmen Consoli - Confusa E Felice ---
$ws.Range("A2").Value = "Carmen Consoli - Confusa E Felice CD"
$ws.Range("B2").Value = "Carmen Consoli - Confusa E Felice`nCD, Album,`nItaly, `nCat. No:`nBarcode: None"
$ws.Range("C2").Value = "07314 537 179-2 01 / 51292784"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "453"
$ws.Range("G2").Value = "5"
$ws.Range("H2").Value = $shipNote
$ws.Range("I2").Value = "ALTA"
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = "hoy"
$ws.Range("L2").Value = "Otros"
$ws.Range("M2").Value = "4,5"

# --- Row 3: Cómplices - Preguntas Y Flores ---
$ws.Range("A3").Value = "Cómplices - Preguntas Y Flores CD"
$ws.Range("B3").Value = "Cómplices - Preguntas Y Flores`nCD, Album,`nSpain, 1993`nCat. No:`nBarcode: 743211738828"
$ws.Range("C3").Value = "7 43211 73882 8"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "453"
$ws.Range("G3").Value = "5"
$ws.Range("H3").Value = $shipNote
$ws.Range("I3").Value = "ALTA"
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = "hoy"
$ws.Range("L3").Value = "Otros"
$ws.Range("M3").Value = "4,5"

# --- Row 4: Phil Collins - Dance Into The Light ---
$ws.Range("A4").Value = "Phil Collins - Dance Into The Light CD"
$ws.Range("B4").Value = "Phil Collins - Dance Into The Light`nCD, Album,`nEurope, 1996-10-14`nCat. No:`nBarcode: 706301600023"
$ws.Range("C4").Value = "706301600023"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "453"
$ws.Range("G4").Value = "5"
$ws.Range("H4").Value = $shipNote
$ws.Range("I4").Value = "ALTA"
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = "hoy"
$ws.Range("L4").Value = "Otros"
$ws.Range("M4").Value = "4,5"
